# Apply updated Betfair Back/Lay odds values for 2025-11-11, rows 2-7
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.68
$ws.Range("G2").Value = 1.87
$ws.Range("H2").Value = 4.9
$ws.Range("K2").Value = 4.1
$ws.Range("L2").Value = 1.37
$ws.Range("O2").Value = 1.31
$ws.Range("P2").Value = 1.83
$ws.Range("R2").Value = 1.4
$ws.Range("V2").Value = 1.19
$ws.Range("W2").Value = 2.14
$ws.Range("Z2").Value = 980
# Row 3
$ws.Range("F3").Value = 1.59
$ws.Range("G3").Value = 1.62
$ws.Range("H3").Value = 6.8
$ws.Range("J3").Value = 4.2
$ws.Range("K3").Value = 4.5
$ws.Range("P3").Value = 2.1
$ws.Range("S3").Value = 2.92
$ws.Range("U3").Value = 2.02
$ws.Range("W3").Value = 2.6
$ws.Range("X3").Value = 19.5
$ws.Range("Y3").Value = 24
$ws.Range("AA3").Value = 210
$ws.Range("AB3").Value = 9.800000000000001
$ws.Range("AC3").Value = 10
$ws.Range("AD3").Value = 25
$ws.Range("AF3").Value = 10.5
$ws.Range("AG3").Value = 11
$ws.Range("AI3").Value = 85
$ws.Range("AJ3").Value = 15
$ws.Range("AK3").Value = 17
$ws.Range("AL3").Value = 36
$ws.Range("AM3").Value = 130
$ws.Range("AO3").Value = 120
# Row 4
$ws.Range("F4").Value = 2.04
$ws.Range("G4").Value = 2.1
$ws.Range("H4").Value = 3.55
$ws.Range("J4").Value = 3.85
$ws.Range("K4").Value = 4.1
$ws.Range("R4").Value = 1.47
$ws.Range("S4").Value = 2.78
$ws.Range("U4").Value = 2.28
$ws.Range("W4").Value = 1.9
$ws.Range("X4").Value = 21
$ws.Range("Y4").Value = 20
$ws.Range("AB4").Value = 13.5
$ws.Range("AE4").Value = 44
$ws.Range("AF4").Value = 17
$ws.Range("AH4").Value = 17.5
$ws.Range("AI4").Value = 48
$ws.Range("AM4").Value = 200
$ws.Range("AN4").Value = 13
$ws.Range("AO4").Value = 38
# Row 5
$ws.Range("F5").Value = 2.28
$ws.Range("G5").Value = 2.46
$ws.Range("H5").Value = 3.45
$ws.Range("I5").Value = 4
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 3.3
$ws.Range("N5").Value = 2.78
$ws.Range("O5").Value = 1.45
$ws.Range("P5").Value = 1.59
$ws.Range("Q5").Value = 2.34
$ws.Range("R5").Value = 1.22
$ws.Range("S5").Value = 4.7
$ws.Range("T5").Value = 1.96
$ws.Range("U5").Value = 1.83
$ws.Range("V5").Value = 1.35
$ws.Range("W5").Value = 1.68
$ws.Range("X5").Value = 10
$ws.Range("AB5").Value = 9.6
$ws.Range("AC5").Value = 7.4
$ws.Range("AF5").Value = 14
$ws.Range("AG5").Value = 12.5
$ws.Range("AH5").Value = 24
$ws.Range("AJ5").Value = 38
$ws.Range("AK5").Value = 34
$ws.Range("AM5").Value = 190
# Row 6
$ws.Range("F6").Value = 2.5
$ws.Range("G6").Value = 2.68
$ws.Range("H6").Value = 3.15
$ws.Range("J6").Value = 2.86
$ws.Range("K6").Value = 3.35
$ws.Range("L6").Value = 1.49
$ws.Range("N6").Value = 2.5
$ws.Range("O6").Value = 1.53
$ws.Range("P6").Value = 1.5
$ws.Range("Q6").Value = 2.58
$ws.Range("S6").Value = 5.5
$ws.Range("T6").Value = 2.06
$ws.Range("U6").Value = 1.77
$ws.Range("W6").Value = 1.59
$ws.Range("X6").Value = 90
$ws.Range("AC6").Value = 14
$ws.Range("AF6").Value = 1000
# Row 7
$ws.Range("G7").Value = 4.4
$ws.Range("J7").Value = 3.05
